# Append/refresh scrape: 2025-09-24 06:26 JST
# Replaces the data rows (2-15) of the "ランサーズ" sheet with a fresh
# 7-row result set (rows 2-8), dropping the now-stale rows 9-15 entirely,
# and tweaks a couple of column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- new result set (already sorted by priority score, descending) ---
$rows = @(
    @{ A = "2025-09-24 06:25:59"; B = "【AI活用】データ分析Webサービス開発パートナー募集"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399092"; G = 368; H = "🔥AI,Ai ◆開発" },
    @{ A = "2025-09-24 06:25:59"; B = "Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399398"; G = 298; H = "🔥Python ◆開発,スクレイピング" },
    @{ A = "2025-09-24 06:25:59"; B = "Googleフォーム × スプレッドシート × GAS 自動化(ストレスチェック診断/台帳保存あり)"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399200"; G = 88; H = "◆自動化" },
    @{ A = "2025-09-24 06:25:59"; B = "完全在宅GASエンジニア募集/課題テストからご依頼/時給1,163円~業務フロー効率化をお任せします"; C = "システム開発"; D = "~ 5,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399071"; G = 70; H = "◆効率化" },
    @{ A = "2025-09-24 06:25:59"; B = "【ペットのアバター化】Pawsitiveプロトタイプ開発の依頼"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399313"; G = 68; H = "◆開発" },
    @{ A = "2025-09-24 06:25:59"; B = "【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398932"; G = 18; H = $null },
    @{ A = "2025-09-24 06:25:59"; B = "限定公開 PR 限定公開の仕事"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5399347"; G = 13; H = $null }
)

# Drop every existing hyperlink first (they reference rows that are about
# to be overwritten/removed, and stale rows must not leave orphan links).
$ws.Cells.Hyperlinks.Delete()

# Remove the old rows 9-15 entirely - the refreshed result set only has
# 7 data rows (2-8) this time.
$ws.Range("A9:H15").EntireRow.Delete()

# Write the refreshed data into rows 2-8.
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G

    if ($null -eq $data.H) {
        $ws.Cells.Item($r, 8).ClearContents()
    } else {
        $ws.Cells.Item($r, 8).Value = $data.H
    }

    $linkCell = $ws.Cells.Item($r, 6)
    $ws.Hyperlinks.Add($linkCell, $data.F) | Out-Null
    $linkCell.Style = "Hyperlink"
}

# Column width tweaks.
# NOTE: the xlsx-exported <col width="..."> value = ColumnWidth + 5/6, so
# back the 5/6 padding out here to land on the exact target stored widths
# of 52 (col B) and 21 (col H).
$ws.Columns.Item(2).ColumnWidth = 52 - (5 / 6)
$ws.Columns.Item(8).ColumnWidth = 21 - (5 / 6)
